$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting all existing data
# (Equipos/Puntos/etc. table) one column to the right.
$ws.Columns.Item(1).Insert()

# Copy the header style (bold + border + centered/top alignment) from the
# shifted header row onto the new rank column so it matches the rest of
# row 1's formatting, then fill in the rank numbers 1-6 for each team.
$ws.Range("B1").Copy()
$ws.Range("A2:A7").PasteSpecial(-4122)

$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5
$ws.Range("A7").Value = 6
